# Notificación - Muestreo Datos.xlsx
# "Creacion de documento Maestro de muestreo de datos, arreglo de muestreos
#  en centro de evento y pagos, y arreglo de diagramas"
#
# On the "Notificación" sheet, the old IdServicio/IdPago sample columns
# (B and C) are replaced by a single "Nombre" column, and the old
# concatenation formula (A-B-C) is replaced with a simpler "A B" formula
# that references the new Nombre column. The now-unused IdServicio column
# is removed entirely (columns shift left).

$wb = $excel.ActiveWorkbook
# 3rd sheet = "Notificación" (accessed by index to avoid any encoding
# issues with the accented sheet name).
$ws = $wb.Worksheets.Item(3)

# Remove the old column C (IdServicio); column D (Combinacion única)
# shifts left to become the new column C.
$ws.Columns.Item(3).Delete()

# Replace the old numeric "IdServicio" sample data in column B with the
# new "Nombre" sample text values (set the data rows first, then the
# header, matching the original authoring order).
$ws.Range("B2").Value2 = "Reserva de Cita "
$ws.Range("B1").Value2 = "Nombre"
$ws.Range("B3").Value2 = "Oferta Shampoo"

# Update the "Combinacion única" formulas (now in column C) to just join
# the identifier and the new name with a space instead of the old
# dash-separated A-B-C concatenation.
$ws.Range("C2").Formula = "=A2&"" ""&B2"
$ws.Range("C3").Formula = "=A3&"" ""&B3"

# Update the last active selection on this sheet.
$ws.Range("F6").Select()
